$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-06-16 Sunday" "2024-06-17 Monday"

Replace-Text "999÷7=142, 5" "627÷8=78, 3"
Replace-Text "508÷6=84, 4" "813÷5=162, 3"
Replace-Text "186÷3=62, 0" "368÷6=61, 2"
Replace-Text "428÷5=85, 3" "573÷7=81, 6"
Replace-Text "323÷7=46, 1" "679÷2=339, 1"

Replace-Text "473÷8=59, 1" "419÷9=46, 5"
Replace-Text "670÷6=111, 4" "761÷3=253, 2"
Replace-Text "261÷8=32, 5" "517÷3=172, 1"
Replace-Text "323÷8=40, 3" "433÷8=54, 1"
Replace-Text "116÷5=23, 1" "923÷6=153, 5"

Replace-Text "636÷2=318, 0" "508÷9=56, 4"
Replace-Text "566÷4=141, 2" "774÷7=110, 4"
Replace-Text "837÷2=418, 1" "477÷4=119, 1"
Replace-Text "415÷8=51, 7" "646÷6=107, 4"
Replace-Text "294÷9=32, 6" "671÷9=74, 5"

Replace-Text "132÷2=66, 0" "462÷3=154, 0"
Replace-Text "307÷4=76, 3" "962÷2=481, 0"
Replace-Text "728÷5=145, 3" "353÷6=58, 5"
Replace-Text "306÷3=102, 0" "225÷5=45, 0"
Replace-Text "176÷2=88, 0" "703÷2=351, 1"

Replace-Text "625÷2=312, 1" "400÷7=57, 1"
Replace-Text "644÷7=92, 0" "502÷2=251, 0"
Replace-Text "805÷8=100, 5" "214÷4=53, 2"
Replace-Text "558÷8=69, 6" "147÷4=36, 3"
Replace-Text "902÷2=451, 0" "703÷9=78, 1"
